$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each changed coin row.
# D-column values that look numeric must be forced to Text format first so
# Excel stores the exact original string (e.g. trailing zeros, "1.000") rather
# than silently converting it to a floating point number.

$ws.Range('D2').Value = '28.185.46'
$ws.Range('E2').Value = '  +5.41%  '
$ws.Range('D3').Value = '1.782.70'
$ws.Range('E3').Value = '  +2.92%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '244.54'
$ws.Range('E5').Value = '  +0.86%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9998'
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4909'
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2679'
$ws.Range('E8').Value = '  +2.08%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06272'
$ws.Range('E9').Value = '  +0.80%  '
$ws.Range('D10').Value = '1.779.21'
$ws.Range('E10').Value = '  +2.76%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '16.50'
$ws.Range('E11').Value = '  +3.98%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07030'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.6286'
$ws.Range('E13').Value = '  +2.63%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.661'
$ws.Range('E14').Value = '  +3.59%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '79.99'
$ws.Range('E15').Value = '  +3.39%  '
$ws.Range('D16').Value = '28.156.30'
$ws.Range('E16').Value = '  +6.14%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.000'
$ws.Range('E17').Value = '  +0.28%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.9993'
$ws.Range('E18').Value = '  +0.17%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007242'
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('E20').Value = '  +5.37%  '
$ws.Range('D21').Value = '2.009.48'
$ws.Range('E21').Value = '  +3.09%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.559'
$ws.Range('E22').Value = '  +1.50%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.738'
$ws.Range('E23').Value = '  +1.94%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.259'
$ws.Range('E24').Value = '  +3.01%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '141.06'
$ws.Range('E25').Value = '  +2.18%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '15.77'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.861'
$ws.Range('E27').Value = '  +4.75%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '109.54'
$ws.Range('E28').Value = '  +2.64%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.384'
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.196'
$ws.Range('E30').Value = '  +6.74%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08269'
$ws.Range('E31').Value = '  +3.46%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.768'
$ws.Range('E32').Value = '  +2.48%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.04895'
$ws.Range('E33').Value = '  +9.20%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.075'
$ws.Range('E34').Value = '  +7.25%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.618'
$ws.Range('E35').Value = '  +0.35%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.6516'
$ws.Range('E36').Value = '  +4.42%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.9520'
$ws.Range('E37').Value = '  +1.82%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.601'
$ws.Range('E38').Value = '  +7.48%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.047'
$ws.Range('E39').Value = '  -0.34%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.892'
$ws.Range('E40').Value = '  +4.94%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.01551'
$ws.Range('E41').Value = '  +2.43%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.9996'
$ws.Range('E42').Value = '  +0.23%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '99.87'
$ws.Range('E43').Value = '  +0.49%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.3985'
$ws.Range('E44').Value = '  +3.26%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '7.190'
$ws.Range('E45').Value = '  +4.05%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.1218'
$ws.Range('E46').Value = '  +4.84%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.05437'
$ws.Range('E47').Value = '  +1.02%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '8.037'
$ws.Range('E48').Value = '  +1.66%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.298'
$ws.Range('E49').Value = '  +4.95%  '
$ws.Range('E50').Value = '  +1.55%  '
$ws.Range('E51').Value = '  +2.17%  '
